# Timesheet update: add new entries for week of 2018-01-07 .. 2018-01-15,
# plus a new "Powerpoint" entry, matching commit "update + add powerpoint".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planning")

# --- 1. Insert 10 new rows right above the old "Totaal:" block (was row 64) ---
$ws.Rows("64:73").Insert()

# Copy the formatting of the last existing data row (now row 63) down into
# the freshly inserted rows so borders / fonts / number formats match.
$ws.Range("A63:C63").Copy()
$ws.Range("A64:C73").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 2. Fill in the new data rows ---
# Row 64
$ws.Range("A64").Value2 = 43107
$ws.Range("B64").Value2 = "dossier, disabled days, info page"
$ws.Range("C64").Value2 = 4

# Row 65
$ws.Range("A65").Value2 = 43108
$ws.Range("B65").Value2 = "geen review eigen voertuig, contactform"
$ws.Range("C65").Value2 = 1.5

# Row 66 - blank spacer row (styled, but empty)

# Row 67
$ws.Range("A67").Value2 = 43109
$ws.Range("B67").Value2 = "info-page, message naar beheer, melding nieuwe request, readme"
$ws.Range("C67").Value2 = 4.25

# Row 68
$ws.Range("A68").Value2 = 43110
$ws.Range("B68").Value2 = "Date home page, book disabled days, profile edit na register, change mail, bugs fixen, change pass"
$ws.Range("C68").Value2 = 4

# Row 69
$ws.Range("A69").Value2 = 43111
$ws.Range("B69").Value2 = "add images"
$ws.Range("C69").Value2 = 3.5

# Row 70
$ws.Range("A70").Value2 = 43112
$ws.Range("B70").Value2 = "componentents detail"
$ws.Range("C70").Value2 = 2.5

# Row 71
$ws.Range("A71").Value2 = 43113
$ws.Range("B71").Value2 = "detail + afwerken, request en bookings"
$ws.Range("C71").Value2 = 3.75

# Row 72
$ws.Range("A72").Value2 = 43114
$ws.Range("B72").Value2 = "afwerken + detail"
$ws.Range("C72").Value2 = 1.5

# Row 73
$ws.Range("A73").Value2 = 43115
$ws.Range("B73").Value2 = "Powerpoint"
$ws.Range("C73").Value2 = 1.25

# Set row heights: default 21, taller (41) for the two rows whose text wraps
# onto a second line.
$ws.Rows("64:73").RowHeight = 21
$ws.Rows("67:68").RowHeight = 41

# --- 3. Old "row 73..82" filler block is removed; row 83 loses its B/C cells ---
$ws.Range("B83").ClearContents()
$ws.Range("C83").ClearContents()

# --- 4. Ten new blank filler rows appended after the old last filler row (96) ---
$ws.Range("A90:C96").Copy()
$ws.Range("A97").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
For ($r = 97; $r -le 106; $r++) {
    $ws.Cells.Item($r, 1).ClearContents()
    $ws.Cells.Item($r, 2).ClearContents()
    $ws.Cells.Item($r, 3).ClearContents()
}

# --- 5. Tail spacing rows (blank, just custom heights) shift from 140-145 to 150-155 ---
$ws.Rows("150").RowHeight = 46
$ws.Rows("151").RowHeight = 17
$ws.Rows("152").RowHeight = 17
$ws.Rows("153").RowHeight = 17
$ws.Rows("154").RowHeight = 17
$ws.Rows("155").RowHeight = 17

# --- 6. Selection / view bookkeeping ---
$ws.Range("E72").Select()

$win = $wb.Windows.Item(1)
$win.Left = 4660
$win.Top = 720
$win.Width = 23160
$win.Height = 15660
